$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the worksheet ---
$ws.Name = "Profile"

# --- tweak a few header labels ---
$ws.Range("F1").Value = "Project Start Date "
$ws.Range("S1").Value = "Actual Spend"
$ws.Range("T1").Value = "Total Revenue recognized"

# --- column widths (A:T) ---
$ws.Range("A1:T1").EntireColumn.ColumnWidth = 9.1666666667

# --- new data rows (2-8) ---

# Row 2
$ws.Range("A2").Value = "fixedPrice/2023-08-24T08:30:35.773Z/client2"
$ws.Range("B2").Value = "project element"
$ws.Range("C2").Value = "'"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "fixedPrice"
$ws.Range("E2").Value = "inProgress"
$ws.Range("F2").Value = "2023-08-24T08:30:35.773Z"
$ws.Range("G2").Value = "Ahmed Shalaab"
$ws.Range("H2").Value = "client2pm"
$ws.Range("K2").Value = "USD"
$ws.Range("L2").Value = 35234
$ws.Range("M2").Value = "Signed"
$ws.Range("N2").Value = "referenceNumber1"
$ws.Range("O2").Value = 12345
$ws.Range("Q2").Value = "client2"

# Row 3
$ws.Range("A3").Value = "TnM/2023-08-24T10:33:33.752Z/client1"
$ws.Range("B3").Value = "iprojectnase2"
$ws.Range("C3").Value = "'"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "TnM"
$ws.Range("E3").Value = "onHold"
$ws.Range("F3").Value = "2023-08-24T10:33:33.752Z"
$ws.Range("G3").Value = "obaid saafan"
$ws.Range("H3").Value = "client1pm"
$ws.Range("K3").Value = "USD"
$ws.Range("L3").Value = 35234
$ws.Range("M3").Value = "Signed"
$ws.Range("N3").Value = "referenceNumber2"
$ws.Range("O3").Value = 12345
$ws.Range("Q3").Value = "client1"

# Row 4
$ws.Range("A4").Value = "TnM/2023-09-16T17:30:36.885Z/client1"
$ws.Range("B4").Value = "iprojectndse2"
$ws.Range("C4").Value = "'"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "TnM"
$ws.Range("E4").Value = "notStarted"
$ws.Range("F4").Value = "2023-09-16T17:30:36.885Z"
$ws.Range("G4").Value = "obaid saafan"
$ws.Range("H4").Value = "client1pm"
$ws.Range("K4").Value = "USD"
$ws.Range("L4").Value = 35234
$ws.Range("M4").Value = "Signed"
$ws.Range("N4").Value = "referenceNumber2"
$ws.Range("O4").Value = 12345
$ws.Range("Q4").Value = "client1"

# Row 5
$ws.Range("A5").Value = "fixedPrice/2023-08-24T08:30:01.277Z/client2"
$ws.Range("B5").Value = "monitor"
$ws.Range("C5").Value = "'"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "fixedPrice"
$ws.Range("E5").Value = "cancelled"
$ws.Range("F5").Value = "2023-08-24T08:30:01.277Z"
$ws.Range("G5").Value = "Ahmed Shalaab"
$ws.Range("H5").Value = "client2pm"
$ws.Range("K5").Value = "USD"
$ws.Range("L5").Value = 35234
$ws.Range("M5").Value = "Signed"
$ws.Range("N5").Value = "referenceNumber1"
$ws.Range("O5").Value = 12345
$ws.Range("Q5").Value = "client2"

# Row 6
$ws.Range("A6").Value = "fixedPrice/2023-08-28T08:15:58.241Z/client2"
$ws.Range("B6").Value = "project element new"
$ws.Range("C6").Value = "'"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "fixedPrice"
$ws.Range("E6").Value = "finished"
$ws.Range("F6").Value = "2023-08-28T08:15:58.241Z"
$ws.Range("G6").Value = "Ahmed Shalaab"
$ws.Range("H6").Value = "client2pm"
$ws.Range("K6").Value = "USD"
$ws.Range("L6").Value = 35234
$ws.Range("M6").Value = "Signed"
$ws.Range("N6").Value = "referenceNumber1"
$ws.Range("O6").Value = 12345
$ws.Range("Q6").Value = "client2"

# Row 7
$ws.Range("A7").Value = "fixedPrice/2023-08-24T08:30:26.178Z/client2"
$ws.Range("B7").Value = "newProject"
$ws.Range("C7").Value = "'"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "fixedPrice"
$ws.Range("E7").Value = "notStarted"
$ws.Range("F7").Value = "2023-08-24T08:30:26.178Z"
$ws.Range("G7").Value = "Ahmed Shalaab"
$ws.Range("H7").Value = "client2pm"
$ws.Range("K7").Value = "USD"
$ws.Range("L7").Value = 35234
$ws.Range("M7").Value = "Signed"
$ws.Range("N7").Value = "referenceNumber1"
$ws.Range("O7").Value = 12345
$ws.Range("Q7").Value = "client2"

# Row 8
$ws.Range("A8").Value = "TnM/3122-06-12T11:31:00Z/client1"
$ws.Range("B8").Value = "iprojectnase"
$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "TnM"
$ws.Range("E8").Value = "notStarted"
$ws.Range("F8").Value = "3122-06-12T11:31:00.000Z"
$ws.Range("G8").Value = "obaid saafan"
$ws.Range("H8").Value = "client1pm"
$ws.Range("K8").Value = "AED"
$ws.Range("L8").Value = 35234
$ws.Range("M8").Value = "Signed"
$ws.Range("N8").Value = "referenceNumber1"
$ws.Range("O8").Value = 12345
$ws.Range("P8").Value = 12345
$ws.Range("Q8").Value = "client1"
